# Updated symbol list on Sat Feb  4 08:50:13 UTC 2023 with GitHub Actions
# Refresh the Price (column D) and Volume(1h) (column E) figures for the
# rows whose crypto market snapshot changed. NumberFormat is forced to
# text ("@") before each write so values like "329.00" / "1.71%" keep
# their exact textual representation instead of being re-interpreted as
# a number/percentage.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "329.00"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "1.71%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "41.30"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "4.89%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.639"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-1.65%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08169"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "2.19%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2.015"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "2.37%"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "1.56%"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-0.27%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.949"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "0.08%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9206"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-0.60%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1273"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "1.99%"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "0.57%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09383"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "2.11%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.03801"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "4.59%"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.90%"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-0.86%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.006137"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.85%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.445"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "2.83%"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-1.52%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "8.335"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-4.24%"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "0.03%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.04388"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.59%"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-0.48%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.004339"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-5.42%"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "4.32%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02808"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "11.15%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05410"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "1.02%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007435"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "0.12%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1416"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "0.89%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.008950"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-5.93%"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "2.51%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.01154"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "7.59%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006589"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-2.86%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.01%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.003193"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "7.41%"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-0.48%"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.01%"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.01%"
